# Insert a new weekly price record as row 97, pushing the existing
# records (old rows 97-143) down by one (new rows 98-144).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(97).Insert()

$ws.Range("A97").Value = 11
$ws.Range("B97").Value = "Vega Monumental Concepción"
$ws.Range("C97").Value = "Bíobío"
$ws.Range("D97").Value = 44876
$ws.Range("E97").Value = 8
$ws.Range("F97").Value = 100112021
$ws.Range("G97").Value = "Ají"
$ws.Range("H97").Value = "Inferno"
$ws.Range("I97").Value = "Primera"
$ws.Range("J97").Value = 80
$ws.Range("K97").Value = 16000
$ws.Range("L97").Value = 17000
$ws.Range("M97").Value = 16625
$ws.Range("N97").Value = "`$/caja 10 kilos"
$ws.Range("O97").Value = "Región de Arica y Parinacota"
$ws.Range("P97").Value = 1662
$ws.Range("Q97").Value = 10
$ws.Range("R97").Value = "Hortaliza"
